$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain numeric-looking price strings to remain text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated price / volume values
$ws.Range("D2").Value = '26.911.66'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '1.553.04'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  +0.61%  '
$ws.Range("D5").Value = '206.73'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("D8").Value = '21.71'
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = '1.773.86'
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").Value = '1.554.10'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("E14").Value = '  +1.54%  '
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = '26.911.63'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").Value = '61.65'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '216.97'
$ws.Range("E18").Value = '  +2.36%  '
$ws.Range("D19").Value = '0.0₃0688'
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").Value = '9.24'
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").Value = '153.76'
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  +1.15%  '
$ws.Range("E30").Value = '  +2.91%  '
$ws.Range("D31").Value = '1.10'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").Value = '1.420.07'
$ws.Range("E33").Value = '  +4.36%  '
$ws.Range("E34").Value = '  +2.90%  '
$ws.Range("E35").Value = '  +3.86%  '
$ws.Range("D36").Value = '0.960'
$ws.Range("E36").Value = '  +1.68%  '
$ws.Range("E37").Value = '  +0.99%  '
$ws.Range("E38").Value = '  +0.98%  '
$ws.Range("D39").Value = '0.520'
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  +1.37%  '
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("D42").Value = '5.69'
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("D44").Value = '2.27'
$ws.Range("E44").Value = '  +3.25%  '
$ws.Range("E45").Value = '  +1.92%  '
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Value = '1.688.59'
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("D48").Value = '86.17'
$ws.Range("E48").Value = '  +0.99%  '
$ws.Range("D49").Value = '0.0524'
$ws.Range("E49").Value = '  +4.22%  '
$ws.Range("E50").Value = '  +4.63%  '
$ws.Range("E51").Value = '  +1.55%  '

